# "adding monitoring and fix minor bugs"
#
# - PENCACAH (G2) used to hold a surveyor's name; replace it with a
#   monitoring contact e-mail address, clearing the stray number-format /
#   right-alignment that had been applied to that one cell and making sure
#   it has no fill (minor cosmetic bugfix).
# - SUMBER (I2) keeps its existing "susenas" value - untouched.
# - Column G needs to be widened so the longer e-mail address fits.
# - Leave the cursor parked on I6 (last place the author clicked before
#   saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- PENCACAH monitoring contact -----------------------------------------
$g2 = $ws.Range("G2")
$g2.Value = "pcl01@bpssumsel.com"
$g2.ClearFormats()
$g2.Interior.Color = 16777215   # explicit "No Fill" (white, renders transparent)

# --- widen column G so the e-mail address fits ---------------------------
$ws.Columns("G").ColumnWidth = 21.3

# --- restore the cursor position saved with the workbook -----------------
$ws.Range("I6").Select() | Out-Null
